$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = 45354
$ws.Range("E1").NumberFormat = "mm-dd-yy"

$ws.Range("E2").Value = 5
$ws.Range("E3").Value = 5
$ws.Range("E4").Value = 3
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 3

$ws.Range("E1").ColumnWidth = 9.5

$ws.Range("E7").Select()
